$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 661.875
$ws.Range("I12").Value = 549
$ws.Range("K12").Value = 549
$ws.Range("M12").Value = -379
$ws.Range("H18").Value = 1046.5
$ws.Range("I18").Value = 952.2222
$ws.Range("K18").Value = 952.2222
$ws.Range("M18").Value = -668.2222
$ws.Range("H19").Value = 900.4
$ws.Range("I19").Value = 650
$ws.Range("J19").Value = 1067.3334
$ws.Range("K19").Value = 650
$ws.Range("L19").Value = 1067.3334
$ws.Range("M19").Value = -475
$ws.Range("N19").Value = -1417.3334
$ws.Range("H33").Value = 110.4
$ws.Range("I33").Value = 56.666668
$ws.Range("K33").Value = 56.666668
$ws.Range("M33").Value = 172.333332
$ws.Range("H41").Value = 902.9
$ws.Range("J41").Value = 1150
$ws.Range("L41").Value = 1150
$ws.Range("N41").Value = -2030
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3866.3333
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 3866.3333
$ws.Range("M64").Value = ""
$ws.Range("N64").Value = -4362.3333
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3866.3333
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 3866.3333
$ws.Range("M67").Value = ""
$ws.Range("N67").Value = -5582.3333
$ws.Range("H80").Value = 265.1
$ws.Range("I80").Value = 162.25
$ws.Range("J80").Value = 333.66666
$ws.Range("K80").Value = 486.75
$ws.Range("L80").Value = 1000.99998
$ws.Range("M80").Value = 511.25
$ws.Range("N80").Value = -2996.99998
$ws.Range("H83").Value = 265.1
$ws.Range("I83").Value = 162.25
$ws.Range("J83").Value = 333.66666
$ws.Range("K83").Value = 1460.25
$ws.Range("L83").Value = 3002.99994
$ws.Range("M83").Value = 3531.75
$ws.Range("N83").Value = -12986.99994
$ws.Range("H98").Value = 849.7143
$ws.Range("I98").Value = 709.8
$ws.Range("J98").Value = 1199.5
$ws.Range("K98").Value = 709.8
$ws.Range("L98").Value = 1199.5
$ws.Range("M98").Value = 788.2
$ws.Range("N98").Value = -4195.5
$ws.Range("H122").Value = 849.7143
$ws.Range("I122").Value = 709.8
$ws.Range("J122").Value = 1199.5
$ws.Range("K122").Value = 2129.4
$ws.Range("L122").Value = 3598.5
$ws.Range("M122").Value = 320.6000000000004
$ws.Range("N122").Value = -8498.5
$ws.Range("H132").Value = 1736.0385
$ws.Range("I132").Value = 1285.48
$ws.Range("K132").Value = 3856.44
$ws.Range("M132").Value = -1326.44

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7400.5713
$ws.Range("I2").Value = 1061.2
$ws.Range("J2").Value = 23249
$ws.Range("K2").Value = 1061.2
$ws.Range("L2").Value = 23249
$ws.Range("M2").Value = -948.2
$ws.Range("N2").Value = -23475
$ws.Range("H110").Value = 111113350
$ws.Range("I110").Value = 142859180
$ws.Range("J110").Value = 2956.5
$ws.Range("K110").Value = 142859180
$ws.Range("L110").Value = 2956.5
$ws.Range("M110").Value = -142857135
$ws.Range("N110").Value = -7046.5
$ws.Range("H116").Value = 7400.5713
$ws.Range("I116").Value = 1061.2
$ws.Range("J116").Value = 23249
$ws.Range("K116").Value = 1061.2
$ws.Range("L116").Value = 23249
$ws.Range("M116").Value = 1232.8
$ws.Range("N116").Value = -27837
$ws.Range("H132").Value = 1511.9048
$ws.Range("J132").Value = 2899.75
$ws.Range("L132").Value = 8699.25
$ws.Range("N132").Value = -13759.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7400.5713
$ws.Range("I3").Value = 1061.2
$ws.Range("J3").Value = 23249
$ws.Range("K3").Value = 1061.2
$ws.Range("L3").Value = 23249
$ws.Range("M3").Value = -947.2
$ws.Range("N3").Value = -23477
$ws.Range("H26").Value = 26537.75
$ws.Range("I26").Value = 26537.75
$ws.Range("K26").Value = 26537.75
$ws.Range("M26").Value = -26245.75
$ws.Range("H86").Value = 2327.8333
$ws.Range("I86").Value = 1244.7142
$ws.Range("J86").Value = 3844.2
$ws.Range("K86").Value = 1244.7142
$ws.Range("L86").Value = 3844.2
$ws.Range("M86").Value = -121.7141999999999
$ws.Range("N86").Value = -6090.2
$ws.Range("H89").Value = 2327.8333
$ws.Range("I89").Value = 1244.7142
$ws.Range("J89").Value = 3844.2
$ws.Range("K89").Value = 6223.571
$ws.Range("L89").Value = 19221
$ws.Range("M89").Value = -607.5709999999999
$ws.Range("N89").Value = -30453
$ws.Range("H94").Value = 371.7143
$ws.Range("I94").Value = 445.4
$ws.Range("J94").Value = 187.5
$ws.Range("K94").Value = 445.4
$ws.Range("L94").Value = 187.5
$ws.Range("M94").Value = 5.600000000000023
$ws.Range("N94").Value = -1089.5
$ws.Range("H96").Value = 17648.25
$ws.Range("I96").Value = 17648.25
$ws.Range("K96").Value = 17648.25
$ws.Range("M96").Value = -14902.25
$ws.Range("H134").Value = 1425.5
$ws.Range("I134").Value = 1030.6428
$ws.Range("J134").Value = 4189.5
$ws.Range("K134").Value = 3091.9284
$ws.Range("L134").Value = 12568.5
$ws.Range("M134").Value = -556.9284000000002
$ws.Range("N134").Value = -17638.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2294.5
$ws.Range("I16").Value = 2294.5
$ws.Range("K16").Value = 2294.5
$ws.Range("M16").Value = -2007.5
$ws.Range("H22").Value = 742.7692
$ws.Range("J22").Value = 2361
$ws.Range("L22").Value = 2361
$ws.Range("N22").Value = -3061
$ws.Range("H31").Value = 1473.9
$ws.Range("I31").Value = 1351.6471
$ws.Range("J31").Value = 2166.6667
$ws.Range("K31").Value = 1351.6471
$ws.Range("L31").Value = 2166.6667
$ws.Range("M31").Value = -1056.6471
$ws.Range("N31").Value = -2756.6667
$ws.Range("H34").Value = 1473.9
$ws.Range("I34").Value = 1351.6471
$ws.Range("J34").Value = 2166.6667
$ws.Range("K34").Value = 1351.6471
$ws.Range("L34").Value = 2166.6667
$ws.Range("M34").Value = -1149.6471
$ws.Range("N34").Value = -2570.6667
$ws.Range("H58").Value = 3363.4375
$ws.Range("I58").Value = 1593.8182
$ws.Range("K58").Value = 1593.8182
$ws.Range("M58").Value = -1390.8182
$ws.Range("H113").Value = 2294.5
$ws.Range("I113").Value = 2294.5
$ws.Range("K113").Value = 2294.5
$ws.Range("M113").Value = -124.5
$ws.Range("H134").Value = 2553
$ws.Range("I134").Value = 2055.6
$ws.Range("K134").Value = 6166.799999999999
$ws.Range("M134").Value = -3631.799999999999
$ws.Range("H136").Value = 3363.4375
$ws.Range("I136").Value = 1593.8182
$ws.Range("K136").Value = 4781.4546
$ws.Range("M136").Value = -2231.4546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 5905
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").Value = ""
$ws.Range("H92").Value = 887.6429000000001
$ws.Range("I92").Value = 599
$ws.Range("J92").Value = 1407.2
$ws.Range("K92").Value = 1797
$ws.Range("L92").Value = 4221.6
$ws.Range("M92").Value = -549
$ws.Range("N92").Value = -6717.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1126.25
$ws.Range("I80").Value = 501.66666
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 501.66666
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = 496.33334
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 1126.25
$ws.Range("I83").Value = 501.66666
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 2508.3333
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = 2483.6667
$ws.Range("N83").Value = -24984
$ws.Range("H132").Value = 112350.11
$ws.Range("I132").Value = 126256.5
$ws.Range("J132").Value = 1099
$ws.Range("K132").Value = 378769.5
$ws.Range("L132").Value = 3297
$ws.Range("M132").Value = -376239.5
$ws.Range("N132").Value = -8357

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6090.6
$ws.Range("I40").Value = 3151.3333
$ws.Range("J40").Value = 10499.5
$ws.Range("K40").Value = 3151.3333
$ws.Range("L40").Value = 10499.5
$ws.Range("M40").Value = -3015.3333
$ws.Range("N40").Value = -10771.5
$ws.Range("H61").Value = 83337176
$ws.Range("I61").Value = 111113570
$ws.Range("J61").Value = 7999.6665
$ws.Range("K61").Value = 111113570
$ws.Range("L61").Value = 7999.6665
$ws.Range("M61").Value = -111113368
$ws.Range("N61").Value = -8403.666499999999
$ws.Range("H68").Value = 4095.7144
$ws.Range("I68").Value = 3111.6667
$ws.Range("K68").Value = 3111.6667
$ws.Range("M68").Value = -2362.6667
$ws.Range("H71").Value = 4095.7144
$ws.Range("I71").Value = 3111.6667
$ws.Range("K71").Value = 15558.3335
$ws.Range("M71").Value = -11814.3335
$ws.Range("H82").Value = 3031.3333
$ws.Range("I82").Value = 527.625
$ws.Range("K82").Value = 527.625
$ws.Range("M82").Value = -166.625
$ws.Range("H85").Value = 3031.3333
$ws.Range("I85").Value = 527.625
$ws.Range("K85").Value = 527.625
$ws.Range("M85").Value = 720.375
$ws.Range("H113").Value = 83337176
$ws.Range("I113").Value = 111113570
$ws.Range("J113").Value = 7999.6665
$ws.Range("K113").Value = 111113570
$ws.Range("L113").Value = 7999.6665
$ws.Range("M113").Value = -111111400
$ws.Range("N113").Value = -12339.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 913.6667
$ws.Range("I81").Value = 470.5
$ws.Range("J81").Value = 1800
$ws.Range("K81").Value = 941
$ws.Range("L81").Value = 3600
$ws.Range("M81").Value = 120
$ws.Range("N81").Value = -5722
$ws.Range("H84").Value = 913.6667
$ws.Range("I84").Value = 470.5
$ws.Range("J84").Value = 1800
$ws.Range("K84").Value = 4705
$ws.Range("L84").Value = 18000
$ws.Range("M84").Value = 599
$ws.Range("N84").Value = -28608
$ws.Range("H96").Value = 1496.6666
$ws.Range("J96").Value = 1495
$ws.Range("L96").Value = 1495
$ws.Range("N96").Value = -4241
$ws.Range("H122").Value = 1915.2667
$ws.Range("I122").Value = 1848.5555
$ws.Range("K122").Value = 5545.666499999999
$ws.Range("M122").Value = -3095.666499999999
$ws.Range("H132").Value = 1552.8889
$ws.Range("I132").Value = 1282.2858
$ws.Range("K132").Value = 3846.8574
$ws.Range("M132").Value = -1316.8574
